$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

# Columns A-D are text fields. "2024-01-24" and "03" look like a date / number
# to Excel's auto-detection, so force them to text via NumberFormat "@" before
# assignment, then clear the format again so the new row matches the plain
# (unstyled) look of the rest of the data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-24"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "16:18:36"

$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "03"
$ws.Cells.Item($row, 4).ClearFormats()

# Columns E-T are numeric resale figures.
$ws.Cells.Item($row, 5).Value = 138512
$ws.Cells.Item($row, 6).Value = 141449
$ws.Cells.Item($row, 7).Value = 171277
$ws.Cells.Item($row, 8).Value = 149067
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123254
$ws.Cells.Item($row, 11).Value = 223803
$ws.Cells.Item($row, 12).Value = 256152
$ws.Cells.Item($row, 13).Value = 184979
$ws.Cells.Item($row, 14).Value = 110026
$ws.Cells.Item($row, 15).Value = 41273
$ws.Cells.Item($row, 16).Value = 30885
$ws.Cells.Item($row, 17).Value = 73419
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42249
$ws.Cells.Item($row, 20).Value = -1
